# Apply the edits described by the diff to df_RSE_settings (Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the RatioConstraints table (rows 11-37) ---

# Rows where only column C changes from 0.15 to 0.5 (column B stays 1)
$onlyCRows = @(11, 12, 16, 17, 18, 22)
foreach ($r in $onlyCRows) {
    $ws.Cells.Item($r, 3).Value = 0.5
}

# Rows where column B changes from 1 to 5 and column C changes from 0.15 to 1
$bAndCRows = @(13, 14, 15, 19, 20, 21, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37)
foreach ($r in $bAndCRows) {
    $ws.Cells.Item($r, 2).Value = 5
    $ws.Cells.Item($r, 3).Value = 1
}

# --- Update the sheet view (scroll position and selection) ---
$ws.Activate()
$excel.ActiveWindow.TopLeftCell = $ws.Range("A27")
$ws.Range("C32").Select()
